# food_system_indicators.xlsx - "first version of data files"
#
# Before tabs: indicator, datapoint, quantity value, named thing
# After  tabs: named thing, indicator, indicator datapoint, indicator datapoint collection, quantity value
#
# NOTE: worksheet handles returned by this host are positional, so they can
# silently start pointing at a different tab once sheets are added/moved/
# renamed around them. To stay safe we always re-look-up a sheet by its
# (stable) name with $wb.Worksheets.Item("...") immediately before using it,
# instead of holding on to a variable across a structural change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "datapoint" -> "indicator datapoint": relabel + reorder header row.
#    Old: A=id  B="datapoint of"  C="has unit"      D="has numeric value"  E=iri F=name G=description
#    New: A="measurement of" B="has unit" C="has numeric value" D=id       E=iri F=name G=description
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("datapoint").Name = "indicator datapoint"

$wsIndicatorDatapoint = $wb.Worksheets.Item("indicator datapoint")
$wsIndicatorDatapoint.Range("A1").Value = "measurement of"
$wsIndicatorDatapoint.Range("B1").Value = "has unit"
$wsIndicatorDatapoint.Range("C1").Value = "has numeric value"
$wsIndicatorDatapoint.Range("D1").Value = "id"
# E1/F1/G1 (iri / name / description) stay as-is.

# ---------------------------------------------------------------------------
# 2. "indicator": reorder header row and move the spatial-scope list
#    validation from column D to column A.
#    Old: A=id  B=name  C=description  D="spatial scope"  E=iri
#    New: A="spatial scope"  B=id  C=iri  D=name  E=description
# ---------------------------------------------------------------------------
$wsIndicator = $wb.Worksheets.Item("indicator")
$wsIndicator.Range("D2:D1048576").Validation.Delete()
$wsIndicator.Range("A1").Value = "spatial scope"
$wsIndicator.Range("B1").Value = "id"
$wsIndicator.Range("C1").Value = "iri"
$wsIndicator.Range("D1").Value = "name"
$wsIndicator.Range("E1").Value = "description"
$wsIndicator.Range("A2:A1048576").Validation.Add(3, 1, 1, '"Eu,EuMemberStates,Regional,Local"')

# ---------------------------------------------------------------------------
# 3. "named thing" and "quantity value" header rows are unchanged; only
#    their tab position moves (handled in step 5).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 4. Brand-new "indicator datapoint collection" sheet, inserted right before
#    "quantity value".
# ---------------------------------------------------------------------------
$wsNewCollection = $wb.Worksheets.Add($wb.Worksheets.Item("quantity value"))
$wsNewCollection.Name = "indicator datapoint collection"
$wsNewCollection.Range("A1").Value = "indicator datapoints"

# ---------------------------------------------------------------------------
# 5. Put every tab into its final order:
#    named thing, indicator, indicator datapoint, indicator datapoint collection, quantity value
#    Each lookup is re-resolved by name right before it is used, since sheet
#    handles are positional and would otherwise go stale after each Move.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("named thing").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("indicator").Move($wb.Worksheets.Item("indicator datapoint"))
$wb.Worksheets.Item("indicator datapoint").Move($wb.Worksheets.Item("indicator datapoint collection"))
$wb.Worksheets.Item("indicator datapoint collection").Move($wb.Worksheets.Item("quantity value"))
